$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 758
$ws.Range("I2").Value = 758
$ws.Range("K2").Value = 758
$ws.Range("M2").Value = -645
$ws.Range("H11").Value = 7.769231
$ws.Range("I11").Value = 7.769231
$ws.Range("K11").Value = 7.769231
$ws.Range("M11").Value = 132.230769
$ws.Range("H74").Value = 4849
$ws.Range("I74").Value = 4205.125
$ws.Range("K74").Value = 4205.125
$ws.Range("M74").Value = -3269.125
$ws.Range("H77").Value = 4849
$ws.Range("I77").Value = 4205.125
$ws.Range("K77").Value = 21025.625
$ws.Range("M77").Value = -16345.625
$ws.Range("H88").Value = 3018.8
$ws.Range("J88").Value = 3041.5
$ws.Range("L88").Value = 3041.5
$ws.Range("N88").Value = -3853.5
$ws.Range("H91").Value = 3018.8
$ws.Range("J91").Value = 3041.5
$ws.Range("L91").Value = 3041.5
$ws.Range("N91").Value = -5849.5
$ws.Range("H103").Value = 427.33334
$ws.Range("I103").Value = 438.9091
$ws.Range("J103").Value = 300
$ws.Range("K103").Value = 1316.7273
$ws.Range("L103").Value = 900
$ws.Range("M103").Value = -730.7273
$ws.Range("N103").Value = -2072
$ws.Range("H132").Value = 4148.846
$ws.Range("I132").Value = 4273.9165
$ws.Range("K132").Value = 12821.7495
$ws.Range("M132").Value = -10291.7495
$ws.Range("H138").Value = 2409.6365
$ws.Range("I138").Value = 805.4286
$ws.Range("J138").Value = 5217
$ws.Range("K138").Value = 2416.2858
$ws.Range("L138").Value = 15651
$ws.Range("M138").Value = 2723.7142
$ws.Range("N138").Value = -25931
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H31").Value = 5499.75
$ws.Range("I31").Value = 5499.75
$ws.Range("K31").Value = 5499.75
$ws.Range("M31").Value = -5205.75
$ws.Range("H122").Value = 1788.3529
$ws.Range("I122").Value = 1788.3529
$ws.Range("K122").Value = 5365.0587
$ws.Range("M122").Value = -2915.0587
$ws.Range("H125").Value = 93666.664
$ws.Range("J125").Value = 93666.664
$ws.Range("L125").Value = 93666.664
$ws.Range("N125").Value = -103506.664
$ws.Range("H132").Value = 4499
$ws.Range("I132").Value = 2998.5454
$ws.Range("K132").Value = 8995.636200000001
$ws.Range("M132").Value = -6465.636200000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 4028
$ws.Range("J5").Value = 5070.3335
$ws.Range("L5").Value = 5070.3335
$ws.Range("N5").Value = -5296.3335
$ws.Range("H7").Value = 27802.5
$ws.Range("I7").Value = 50
$ws.Range("K7").Value = 50
$ws.Range("M7").Value = 63
$ws.Range("H86").Value = 4206.636
$ws.Range("I86").Value = 4127.3
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4127.3
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3004.3
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4206.636
$ws.Range("I89").Value = 4127.3
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 20636.5
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -15020.5
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 1999.8572
$ws.Range("I99").Value = 1999.8572
$ws.Range("K99").Value = 1999.8572
$ws.Range("M99").Value = -501.8571999999999
$ws.Range("H105").Value = 5141.174
$ws.Range("I105").Value = 4782.8
$ws.Range("K105").Value = 4782.8
$ws.Range("M105").Value = -3035.8
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 2799.5454
$ws.Range("I134").Value = 2779.5
$ws.Range("K134").Value = 8338.5
$ws.Range("M134").Value = -5803.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1400
$ws.Range("I2").Value = 1350
$ws.Range("K2").Value = 1350
$ws.Range("M2").Value = -1237
$ws.Range("H6").Value = 1401779.8
$ws.Range("I6").Value = 1751225
$ws.Range("J6").Value = 3999
$ws.Range("K6").Value = 1751225
$ws.Range("L6").Value = 3999
$ws.Range("M6").Value = -1751112
$ws.Range("N6").Value = -4225
$ws.Range("H31").Value = 2247
$ws.Range("I31").Value = 2247
$ws.Range("K31").Value = 2247
$ws.Range("M31").Value = -1952
$ws.Range("H34").Value = 2247
$ws.Range("I34").Value = 2247
$ws.Range("K34").Value = 2247
$ws.Range("M34").Value = -2045
$ws.Range("H58").Value = 5830
$ws.Range("I58").Value = 5830
$ws.Range("K58").Value = 5830
$ws.Range("M58").Value = -5627
$ws.Range("H136").Value = 5830
$ws.Range("I136").Value = 5830
$ws.Range("K136").Value = 17490
$ws.Range("M136").Value = -14940
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 575.8
$ws.Range("I23").Value = 399.5
$ws.Range("K23").Value = 1198.5
$ws.Range("M23").Value = -963.5
$ws.Range("H51").Value = 1999
$ws.Range("I51").Value = 1999
$ws.Range("K51").Value = 5997
$ws.Range("M51").Value = -5537
$ws.Range("H97").Value = 1250
$ws.Range("I97").Value = 1250
$ws.Range("K97").Value = 3750
$ws.Range("M97").Value = -3254
$ws.Range("H122").Value = 502.875
$ws.Range("J122").Value = 379.33334
$ws.Range("L122").Value = 3414.00006
$ws.Range("N122").Value = -8314.00006
$ws.Range("H131").Value = 2140.7778
$ws.Range("J131").Value = 2066.6667
$ws.Range("L131").Value = 6200.000100000001
$ws.Range("N131").Value = -16280.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118.111115
$ws.Range("I2").Value = 134.85715
$ws.Range("J2").Value = 59.5
$ws.Range("K2").Value = 134.85715
$ws.Range("L2").Value = 59.5
$ws.Range("M2").Value = -21.85714999999999
$ws.Range("N2").Value = -285.5
$ws.Range("H102").Value = 2886.0833
$ws.Range("I102").Value = 2466.6365
$ws.Range("K102").Value = 2466.6365
$ws.Range("M102").Value = -844.6365000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1543.8889
$ws.Range("I22").Value = 459
$ws.Range("K22").Value = 459
$ws.Range("M22").Value = -164
$ws.Range("H27").Value = 1543.8889
$ws.Range("I27").Value = 459
$ws.Range("K27").Value = 459
$ws.Range("M27").Value = -352
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -864
$ws.Range("H46").Value = 1744.5454
$ws.Range("I46").Value = 1028.2858
$ws.Range("J46").Value = 2998
$ws.Range("K46").Value = 1028.2858
$ws.Range("L46").Value = 2998
$ws.Range("M46").Value = -840.2858000000001
$ws.Range("N46").Value = -3374
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H122").Value = 4024.818
$ws.Range("I122").Value = 3909.875
$ws.Range("K122").Value = 11729.625
$ws.Range("M122").Value = -9279.625
$ws.Range("H139").Value = 15000
$ws.Range("J139").Value = 15000
$ws.Range("L139").Value = 15000
$ws.Range("N139").Value = -25280
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 350
$ws.Range("J5").Value = 350
$ws.Range("L5").Value = 350
$ws.Range("N5").Value = -574
$ws.Range("H24").Value = 35000
$ws.Range("I24").Value = 35000
$ws.Range("K24").Value = 35000
$ws.Range("M24").Value = -34770
$ws.Range("H81").Value = 4175
$ws.Range("J81").Value = 6133.3335
$ws.Range("L81").Value = 12266.667
$ws.Range("N81").Value = -14388.667
$ws.Range("H84").Value = 4175
$ws.Range("J84").Value = 6133.3335
$ws.Range("L84").Value = 61333.335
$ws.Range("N84").Value = -71941.33499999999
$ws.Range("H96").Value = 1331.3334
$ws.Range("I96").Value = 1331.3334
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1331.3334
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 41.66660000000002
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 3091.4092
$ws.Range("I122").Value = 3069
$ws.Range("J122").Value = 3233.3333
$ws.Range("K122").Value = 9207
$ws.Range("L122").Value = 9699.999899999999
$ws.Range("M122").Value = -6757
$ws.Range("N122").Value = -14599.9999
$ws.Range("H136").Value = 3556.875
$ws.Range("I136").Value = 3566
$ws.Range("K136").Value = 10698
$ws.Range("M136").Value = -8148
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
